$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "BAG_TYPE_EQUIPMENT"
$ws.Range("D6").Value = "BAG_TYPE_ETC"
$ws.Range("D7").Value = "BAG_TYPE_USEABLE"

$ws.Columns.Item(4).ColumnWidth = 19.125

$ws.Range("N14").Select()
